$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Pesquisa"
$h12 = $ws2.Range("H12")
$ws1.Range("B2").Copy()
$h12.PasteSpecial(-4122)
$h12.NumberFormat = "general"
$h12.HorizontalAlignment = 1
$h12.VerticalAlignment = -4107
$h12.Font.Size = 11
